$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) and "Volume(1h)" (E) columns store plain-text values
# (e.g. "1.00", "313.68") rather than numbers, matching the source sheet,
# which keeps trailing zeros / thousand-dot formatting intact.
# Force text format before writing numeric-looking strings to Price cells
# so Excel does not silently convert them to numbers and drop formatting.

# Row 2
$ws.Range("D2").Value = "44.493.73"
$ws.Range("E2").Value = "  +3.57%  "

# Row 3
$ws.Range("D3").Value = "2.422.61"
$ws.Range("E3").Value = "  +2.47%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.68"
$ws.Range("E5").Value = "  +3.69%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.00"
$ws.Range("E6").Value = "  +5.51%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.513"
$ws.Range("E7").Value = "  +2.21%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("E9").Value = "  +4.66%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.25"
$ws.Range("E10").Value = "  +3.32%  "

# Row 11
$ws.Range("B11").Value = "Chainlink"
$ws.Range("C11").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.21"
$ws.Range("E11").Value = "  +4.56%  "

# Row 12
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0798"
$ws.Range("E12").Value = "  +1.70%  "

# Row 13
$ws.Range("E13").Value = "  -2.21%  "

# Row 14
$ws.Range("E14").Value = "  +3.00%  "

# Row 15
$ws.Range("D15").Value = "2.799.68"
$ws.Range("E15").Value = "  +2.53%  "

# Row 16
$ws.Range("D16").Value = "2.419.25"
$ws.Range("E16").Value = "  +1.94%  "

# Row 17
$ws.Range("E17").Value = "  +5.27%  "

# Row 18
$ws.Range("D18").Value = "44.363.78"
$ws.Range("E18").Value = "  +3.36%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.42"
$ws.Range("E19").Value = "  +4.05%  "

# Row 20
$ws.Range("E20").Value = "  +2.04%  "

# Row 21
$ws.Range("E21").Value = "  +3.61%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.62"
$ws.Range("E22").Value = "  +0.82%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.28"
$ws.Range("E23").Value = "  +2.59%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.29"
$ws.Range("E24").Value = "  +5.72%  "

# Row 25
$ws.Range("E25").Value = "  +1.43%  "

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("E27").Value = "  +2.38%  "

# Row 28
$ws.Range("E28").Value = "  -3.79%  "

# Row 29
$ws.Range("E29").Value = "  +2.78%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.26"
$ws.Range("E30").Value = "  +3.98%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.45"
$ws.Range("E31").Value = "  +0.70%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("E32").Value = "  +17.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.30"
$ws.Range("E33").Value = "  +10.31%  "

# Row 34
$ws.Range("E34").Value = "  +3.01%  "

# Row 35
$ws.Range("E35").Value = "  +0.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0769"
$ws.Range("E36").Value = "  +7.54%  "

# Row 37
$ws.Range("E37").Value = "  +2.01%  "

# Row 38
$ws.Range("E38").Value = "  +3.22%  "

# Row 39
$ws.Range("E39").Value = "  +1.61%  "

# Row 40
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.23"
$ws.Range("E40").Value = "  -1.87%  "

# Row 41
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.19"
$ws.Range("E41").Value = "  -7.79%  "

# Row 42
$ws.Range("E42").Value = "  +1.32%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.88"
$ws.Range("E43").Value = "  -1.25%  "

# Row 44
$ws.Range("E44").Value = "  +3.57%  "

# Row 45
$ws.Range("D45").Value = "1.943.03"
$ws.Range("E45").Value = "  +0.52%  "

# Row 47
$ws.Range("E47").Value = "  +8.09%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.48"
$ws.Range("E48").Value = "  +3.28%  "

# Row 49
$ws.Range("E49").Value = "  +9.73%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.00"
$ws.Range("E50").Value = "  +6.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.12"
$ws.Range("E51").Value = "  +3.75%  "
